$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.986.36"
$ws.Range("E2").Value = "  -7.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.64"
$ws.Range("E3").Value = "  -3.23%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.84"
$ws.Range("E5").Value = "  -7.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "120.63"
$ws.Range("E6").Value = "  -7.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.501.71"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  -12.22%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.667"
$ws.Range("E10").Value = "  -12.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("E11").Value = "  -24.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000306"
$ws.Range("E12").Value = "  -28.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.02"
$ws.Range("E13").Value = "  -9.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.084.93"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.96"
$ws.Range("E15").Value = "  -8.94%  "
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.499.88"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.61"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.43"
$ws.Range("E19").Value = "  -8.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.110.62"
$ws.Range("E20").Value = "  -7.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -10.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.90"
$ws.Range("E22").Value = "  -16.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.73"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.03"
$ws.Range("E24").Value = "  -10.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.83"
$ws.Range("E25").Value = "  -7.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.46"
$ws.Range("E26").Value = "  +9.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.05"
$ws.Range("E27").Value = "  -6.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.93"
$ws.Range("E28").Value = "  -10.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.51"
$ws.Range("E29").Value = "  -16.57%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.64"
$ws.Range("E30").Value = "  -5.42%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.58"
$ws.Range("E31").Value = "  -5.28%  "
$ws.Range("E32").Value = "  -8.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.62"
$ws.Range("E33").Value = "  -8.27%  "
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.147"
$ws.Range("E35").Value = "  -7.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "35.95"
$ws.Range("E36").Value = "  -10.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.85"
$ws.Range("E37").Value = "  -5.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0429"
$ws.Range("E38").Value = "  -12.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.63"
$ws.Range("E40").Value = "  +2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.128"
$ws.Range("E41").Value = "  -13.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.66"
$ws.Range("E42").Value = "  +21.55%  "
$ws.Range("E43").Value = "  +14.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "138.40"
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0597"
$ws.Range("E45").Value = "  -25.28%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  -9.67%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.03"
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.01"
$ws.Range("E49").Value = "  -6.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -11.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.269"
$ws.Range("E51").Value = "  -10.62%  "
